$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (format) of the last existing data row (255) down through
# the new rows 256-269 so the new date cells in column A pick up the same
# cellXf (border/font/alignment/number-format) as the existing rows, without
# minting a brand new style.
$ws.Range("A255").Copy() | Out-Null
$ws.Range("A256:A269").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(256, 1).Value = 44330
$ws.Cells.Item(256, 2).Value = 4
$ws.Cells.Item(256, 3).Value = 23
$ws.Cells.Item(256, 4).Value = 66.92661351335623

$ws.Cells.Item(257, 1).Value = 44331
$ws.Cells.Item(257, 2).Value = 3
$ws.Cells.Item(257, 3).Value = 21
$ws.Cells.Item(257, 4).Value = 61.10690799045568

$ws.Cells.Item(258, 1).Value = 44332
$ws.Cells.Item(258, 2).Value = 3
$ws.Cells.Item(258, 3).Value = 22
$ws.Cells.Item(258, 4).Value = 64.01676075190595

$ws.Cells.Item(259, 1).Value = 44333
$ws.Cells.Item(259, 2).Value = 11
$ws.Cells.Item(259, 3).Value = 31
$ws.Cells.Item(259, 4).Value = 90.2054356049584

$ws.Cells.Item(260, 1).Value = 44334
$ws.Cells.Item(260, 2).Value = 2
$ws.Cells.Item(260, 3).Value = 32
$ws.Cells.Item(260, 4).Value = 93.11528836640865

$ws.Cells.Item(261, 1).Value = 44335
$ws.Cells.Item(261, 2).Value = 3
$ws.Cells.Item(261, 3).Value = 31
$ws.Cells.Item(261, 4).Value = 90.2054356049584

$ws.Cells.Item(262, 1).Value = 44336
$ws.Cells.Item(262, 2).Value = 5
$ws.Cells.Item(262, 3).Value = 31
$ws.Cells.Item(262, 4).Value = 90.2054356049584

$ws.Cells.Item(263, 1).Value = 44337
$ws.Cells.Item(263, 2).Value = 1
$ws.Cells.Item(263, 3).Value = 28
$ws.Cells.Item(263, 4).Value = 81.47587732060757

$ws.Cells.Item(264, 1).Value = 44338
$ws.Cells.Item(264, 2).Value = 3
$ws.Cells.Item(264, 3).Value = 28
$ws.Cells.Item(264, 4).Value = 81.47587732060757

$ws.Cells.Item(265, 1).Value = 44339
$ws.Cells.Item(265, 2).Value = 5
$ws.Cells.Item(265, 3).Value = 30
$ws.Cells.Item(265, 4).Value = 87.29558284350811

$ws.Cells.Item(266, 1).Value = 44340
$ws.Cells.Item(266, 2).Value = 0
$ws.Cells.Item(266, 3).Value = 19
$ws.Cells.Item(266, 4).Value = 55.28720246755515

$ws.Cells.Item(267, 1).Value = 44341
$ws.Cells.Item(267, 2).Value = 1
$ws.Cells.Item(267, 3).Value = 18
$ws.Cells.Item(267, 4).Value = 52.37734970610487

$ws.Cells.Item(268, 1).Value = 44342
$ws.Cells.Item(268, 2).Value = 3
$ws.Cells.Item(268, 3).Value = 18
$ws.Cells.Item(268, 4).Value = 52.37734970610487

$ws.Cells.Item(269, 1).Value = 44343
$ws.Cells.Item(269, 2).Value = 6
$ws.Cells.Item(269, 3).Value = 19
$ws.Cells.Item(269, 4).Value = 55.28720246755515
